# Update TPM-derived NATMI LR-pair statistics for Reln-Itgb1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.1779566666666667
$ws.Range("H2").Value = 0.5338700000000001
$ws.Range("I2").Value = 0.01192558037548992
$ws.Range("J2").Value = 0.01192558037548992
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 13.71859460732222
$ws.Range("R2").Value = 123.4673514659
$ws.Range("S2").Value = 0.002866690895047208
$ws.Range("T2").Value = 0.002866690895047207
$ws.Range("G3").Value = 0.1779566666666667
$ws.Range("H3").Value = 0.5338700000000001
$ws.Range("I3").Value = 0.01192558037548992
$ws.Range("J3").Value = 0.01192558037548992
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 18.07684484371556
$ws.Range("R3").Value = 162.69160359344
$ws.Range("S3").Value = 0.003777407818217868
$ws.Range("T3").Value = 0.003777407818217867
$ws.Range("G4").Value = 0.1779566666666667
$ws.Range("H4").Value = 0.5338700000000001
$ws.Range("I4").Value = 0.01192558037548992
$ws.Range("J4").Value = 0.01192558037548992
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 25.27461400712889
$ws.Range("R4").Value = 227.47152606416
$ws.Range("S4").Value = 0.005281481662224848
$ws.Range("T4").Value = 0.005281481662224847
$ws.Range("I5").Value = 0.540575811616083
$ws.Range("J5").Value = 0.540575811616083
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 621.8515309600244
$ws.Range("R5").Value = 5596.663778640221
$ws.Range("S5").Value = 0.1299445149376151
$ws.Range("T5").Value = 0.1299445149376151
$ws.Range("I6").Value = 0.540575811616083
$ws.Range("J6").Value = 0.540575811616083
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.1712264923671837
$ws.Range("T6").Value = 0.1712264923671837
$ws.Range("I7").Value = 0.540575811616083
$ws.Range("J7").Value = 0.540575811616083
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.2394048043112842
$ws.Range("T7").Value = 0.2394048043112842
$ws.Range("G8").Value = 6.677692666666666
$ws.Range("I8").Value = 0.4474986080084269
$ws.Range("J8").Value = 0.4474986080084269
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 514.7801446398289
$ws.Range("R8").Value = 4633.02130175846
$ws.Range("S8").Value = 0.1075704615400201
$ws.Range("T8").Value = 0.1075704615400201
$ws.Range("G9").Value = 6.677692666666666
$ws.Range("I9").Value = 0.4474986080084269
$ws.Range("J9").Value = 0.4474986080084269
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 678.3202703805262
$ws.Range("S9").Value = 0.1417444423926581
$ws.Range("T9").Value = 0.141744442392658
$ws.Range("G10").Value = 6.677692666666666
$ws.Range("I10").Value = 0.4474986080084269
$ws.Range("J10").Value = 0.4474986080084269
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("Q10").Value = 948.4112496014116
$ws.Range("R10").Value = 8535.701246412704
$ws.Range("S10").Value = 0.1981837040757488
$ws.Range("T10").Value = 0.1981837040757488
